$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D to make room for "Status"
$ws.Columns.Item(4).Insert()

# --- Header row (row 1) ---
$header = New-Object 'object[,]' 1,6
$header[0,0] = "Status"
$header[0,1] = "Jan_2026"
$header[0,2] = "Dec_2025"
$header[0,3] = "Oct_2025"
$header[0,4] = "MoM"
$header[0,5] = "QoQ"
$ws.Range("D1:I1").Value = $header

# --- Data rows 2 to 11 (ISIN/Stock/Fund unchanged; insert Status col, shift numeric cols) ---
$data2_11 = New-Object 'object[,]' 10,6
$data2_11[0,0] = "Fresh Entry"
$data2_11[0,1] = 10.101646
$data2_11[0,2] = 0
$data2_11[0,3] = 0
$data2_11[0,4] = 10.101646
$data2_11[0,5] = 10.101646
$data2_11[1,0] = "Adding Consistently"
$data2_11[1,1] = 9.69026
$data2_11[1,2] = 9.588673
$data2_11[1,3] = 9.091605
$data2_11[1,4] = 0.1015870000000003
$data2_11[1,5] = 0.5986550000000008
$data2_11[2,0] = "Adding Consistently"
$data2_11[2,1] = 9.38495
$data2_11[2,2] = 8.925459
$data2_11[2,3] = 8.369671
$data2_11[2,4] = 0.4594909999999999
$data2_11[2,5] = 1.015279
$data2_11[3,0] = "Reducing"
$data2_11[3,1] = 8.55998
$data2_11[3,2] = 8.703604
$data2_11[3,3] = 6.252576
$data2_11[3,4] = -0.1436240000000009
$data2_11[3,5] = 2.307403999999999
$data2_11[4,0] = "Reducing"
$data2_11[4,1] = 7.921467
$data2_11[4,2] = 8.653581
$data2_11[4,3] = 3.75492
$data2_11[4,4] = -0.732114000000001
$data2_11[4,5] = 4.166547
$data2_11[5,0] = "Adding Consistently"
$data2_11[5,1] = 7.343367
$data2_11[5,2] = 6.98031
$data2_11[5,3] = 0
$data2_11[5,4] = 0.3630569999999995
$data2_11[5,5] = 7.343367
$data2_11[6,0] = "Adding Consistently"
$data2_11[6,1] = 7.217659
$data2_11[6,2] = 6.779965
$data2_11[6,3] = 0
$data2_11[6,4] = 0.4376940000000005
$data2_11[6,5] = 7.217659
$data2_11[7,0] = "Fresh Entry"
$data2_11[7,1] = 6.520862
$data2_11[7,2] = 0
$data2_11[7,3] = 0
$data2_11[7,4] = 6.520862
$data2_11[7,5] = 6.520862
$data2_11[8,0] = "Reducing"
$data2_11[8,1] = 2.931496
$data2_11[8,2] = 3.074732
$data2_11[8,3] = 0
$data2_11[8,4] = -0.1432359999999999
$data2_11[8,5] = 2.931496
$data2_11[9,0] = "Adding Consistently"
$data2_11[9,1] = 0.603333
$data2_11[9,2] = 0.534424
$data2_11[9,3] = 0.50197
$data2_11[9,4] = 0.068909
$data2_11[9,5] = 0.101363
$ws.Range("D2:I11").Value = $data2_11

# --- Data rows 12 to 23 (ISIN/Stock Name changed + Status/numeric cols) ---
$data12_23 = New-Object 'object[,]' 12,9
$data12_23[0,0] = "INE263A01024"
$data12_23[0,1] = "Bharat Electronics Ltd"
$data12_23[0,2] = "quant PSU Fund"
$data12_23[0,3] = "Complete Exit"
$data12_23[0,4] = 0
$data12_23[0,5] = 0
$data12_23[0,6] = 8.270709
$data12_23[0,7] = 0
$data12_23[0,8] = -8.270709
$data12_23[1,0] = "INE752E01010"
$data12_23[1,1] = "Power Grid Corporation of India Limited"
$data12_23[1,2] = "quant PSU Fund"
$data12_23[1,3] = "Complete Exit"
$data12_23[1,4] = 0
$data12_23[1,5] = 0
$data12_23[1,6] = 2.752036
$data12_23[1,7] = 0
$data12_23[1,8] = -2.752036
$data12_23[2,0] = "INE510A01028"
$data12_23[2,1] = "Engineers India Limited"
$data12_23[2,2] = "quant PSU Fund"
$data12_23[2,3] = "Complete Exit"
$data12_23[2,4] = 0
$data12_23[2,5] = 1.011007
$data12_23[2,6] = 2.994703
$data12_23[2,7] = -1.011007
$data12_23[2,8] = -2.994703
$data12_23[3,0] = "INE242A01010"
$data12_23[3,1] = "Indian Oil Corp Ltd"
$data12_23[3,2] = "quant PSU Fund"
$data12_23[3,3] = "Complete Exit"
$data12_23[3,4] = 0
$data12_23[3,5] = 0
$data12_23[3,6] = 8.030499
$data12_23[3,7] = 0
$data12_23[3,8] = -8.030499
$data12_23[4,0] = "INE257A01026"
$data12_23[4,1] = "Bharat Heavy Electricals Ltd"
$data12_23[4,2] = "quant PSU Fund"
$data12_23[4,3] = "Complete Exit"
$data12_23[4,4] = 0
$data12_23[4,5] = 8.2214
$data12_23[4,6] = 0
$data12_23[4,7] = -8.2214
$data12_23[4,8] = 0
$data12_23[5,0] = "INE029A01011"
$data12_23[5,1] = "Bharat Petroleum Corp Ltd"
$data12_23[5,2] = "quant PSU Fund"
$data12_23[5,3] = "Complete Exit"
$data12_23[5,4] = 0
$data12_23[5,5] = 0
$data12_23[5,6] = 6.265599
$data12_23[5,7] = 0
$data12_23[5,8] = -6.265599
$data12_23[6,0] = "INE171Z01026"
$data12_23[6,1] = "Bharat Dynamics Limited"
$data12_23[6,2] = "quant PSU Fund"
$data12_23[6,3] = "Complete Exit"
$data12_23[6,4] = 0
$data12_23[6,5] = 3.631098
$data12_23[6,6] = 0
$data12_23[6,7] = -3.631098
$data12_23[6,8] = 0
$data12_23[7,0] = "INE129A01019"
$data12_23[7,1] = "GAIL (India) Limited"
$data12_23[7,2] = "quant PSU Fund"
$data12_23[7,3] = "Complete Exit"
$data12_23[7,4] = 0
$data12_23[7,5] = 0
$data12_23[7,6] = 1.61757
$data12_23[7,7] = 0
$data12_23[7,8] = -1.61757
$data12_23[8,0] = "INE094A01015"
$data12_23[8,1] = "Hindustan Petroleum Corporation Ltd"
$data12_23[8,2] = "quant PSU Fund"
$data12_23[8,3] = "Complete Exit"
$data12_23[8,4] = 0
$data12_23[8,5] = 0
$data12_23[8,6] = 4.348009
$data12_23[8,7] = 0
$data12_23[8,8] = -4.348009
$data12_23[9,0] = "INE062A01020"
$data12_23[9,1] = "State Bank of India"
$data12_23[9,2] = "quant PSU Fund"
$data12_23[9,3] = "Complete Exit"
$data12_23[9,4] = 0
$data12_23[9,5] = 9.260284
$data12_23[9,6] = 8.642809
$data12_23[9,7] = -9.260284
$data12_23[9,8] = -8.642809
$data12_23[10,0] = "INE031A01017"
$data12_23[10,1] = "Housing & Urban Devlopment Company Ltd"
$data12_23[10,2] = "quant PSU Fund"
$data12_23[10,3] = "Complete Exit"
$data12_23[10,4] = 0
$data12_23[10,5] = 6.050591
$data12_23[10,6] = 2.920438
$data12_23[10,7] = -6.050591
$data12_23[10,8] = -2.920438
$data12_23[11,0] = "INE931S01010"
$data12_23[11,1] = "Adani Energy Solutions Limited"
$data12_23[11,2] = "quant PSU Fund"
$data12_23[11,3] = "Complete Exit"
$data12_23[11,4] = 0
$data12_23[11,5] = 0
$data12_23[11,6] = 2.477313
$data12_23[11,7] = 0
$data12_23[11,8] = -2.477313
$ws.Range("A12:I23").Value = $data12_23

Write-Output "done"
